# Auto-generated edit script applying the Pandaemonium_Profits diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) on specific
# rows across the ALC, ARM, CRP, CUL, GSM, LTW, WVR leve-profit tables.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 58000
$ws.Range("J3").Value = 58000
$ws.Range("L3").Value = 58000
$ws.Range("N3").Value = -58228
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").Value = ""
$ws.Range("H76").Value = 4598
$ws.Range("I76").Value = 3805.5557
$ws.Range("K76").Value = 3805.5557
$ws.Range("M76").Value = -3490.5557
$ws.Range("H79").Value = 4598
$ws.Range("I79").Value = 3805.5557
$ws.Range("K79").Value = 3805.5557
$ws.Range("M79").Value = -2713.5557
$ws.Range("H98").Value = 5467.4
$ws.Range("I98").Value = 5400
$ws.Range("J98").Value = 5624.6665
$ws.Range("K98").Value = 5400
$ws.Range("L98").Value = 5624.6665
$ws.Range("M98").Value = -3902
$ws.Range("N98").Value = -8620.666499999999
$ws.Range("H102").Value = 58000
$ws.Range("J102").Value = 58000
$ws.Range("L102").Value = 58000
$ws.Range("N102").Value = -64490
$ws.Range("H107").Value = 694.5833
$ws.Range("I107").Value = 694.5833
$ws.Range("K107").Value = 694.5833
$ws.Range("M107").Value = 1225.4167
$ws.Range("H113").Value = 2738.75
$ws.Range("I113").Value = 2738.75
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2738.75
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 515.25
$ws.Range("N113").Value = ""
$ws.Range("H122").Value = 5467.4
$ws.Range("I122").Value = 5400
$ws.Range("J122").Value = 5624.6665
$ws.Range("K122").Value = 16200
$ws.Range("L122").Value = 16873.9995
$ws.Range("M122").Value = -13750
$ws.Range("N122").Value = -21773.9995
$ws.Range("H127").Value = 1941.9474
$ws.Range("I127").Value = 474.25
$ws.Range("J127").Value = 2333.3333
$ws.Range("K127").Value = 1422.75
$ws.Range("L127").Value = 6999.999899999999
$ws.Range("M127").Value = 3537.25
$ws.Range("N127").Value = -16919.9999
$ws.Range("H137").Value = 1193772.9
$ws.Range("J137").Value = 1963048.9
$ws.Range("L137").Value = 5889146.699999999
$ws.Range("N137").Value = -5894246.699999999
$ws.Range("H138").Value = 3868.6128
$ws.Range("I138").Value = 3241.1428
$ws.Range("J138").Value = 4051.625
$ws.Range("K138").Value = 9723.428400000001
$ws.Range("L138").Value = 12154.875
$ws.Range("M138").Value = -4583.428400000001
$ws.Range("N138").Value = -22434.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 42000
$ws.Range("J104").Value = 42000
$ws.Range("L104").Value = 42000
$ws.Range("N104").Value = -48988

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 21038
$ws.Range("J43").Value = 21038
$ws.Range("L43").Value = 21038
$ws.Range("N43").Value = -21406
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = ""
$ws.Range("N99").Value = ""
$ws.Range("H101").Value = 21038
$ws.Range("J101").Value = 21038
$ws.Range("L101").Value = 21038
$ws.Range("N101").Value = -27528
$ws.Range("H107").Value = 1373.4
$ws.Range("I107").Value = 1373.4
$ws.Range("K107").Value = 1373.4
$ws.Range("M107").Value = 546.5999999999999
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 300
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -188
$ws.Range("N7").Value = ""
$ws.Range("H92").Value = 502
$ws.Range("I92").Value = 502
$ws.Range("K92").Value = 1506
$ws.Range("M92").Value = -258
$ws.Range("H107").Value = 692.0357
$ws.Range("I107").Value = 291.58536
$ws.Range("J107").Value = 1786.6
$ws.Range("K107").Value = 874.7560799999999
$ws.Range("L107").Value = 5359.799999999999
$ws.Range("M107").Value = 1045.24392
$ws.Range("N107").Value = -9199.799999999999
$ws.Range("H108").Value = 1000
$ws.Range("I108").Value = 1000
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 3000
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -120
$ws.Range("N108").Value = ""
$ws.Range("H110").Value = 4632.875
$ws.Range("I110").Value = 3569
$ws.Range("J110").Value = 4784.857
$ws.Range("K110").Value = 10707
$ws.Range("L110").Value = 14354.571
$ws.Range("M110").Value = -6617
$ws.Range("N110").Value = -22534.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""
$ws.Range("H123").Value = 10289.444
$ws.Range("J123").Value = 10289.444
$ws.Range("L123").Value = 10289.444
$ws.Range("N123").Value = -15189.444
$ws.Range("H124").Value = 79800
$ws.Range("J124").Value = 79800
$ws.Range("L124").Value = 79800
$ws.Range("N124").Value = -89620
$ws.Range("H125").Value = 79800
$ws.Range("J125").Value = 79800
$ws.Range("L125").Value = 79800
$ws.Range("N125").Value = -84720
$ws.Range("H126").Value = 3186.6667
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 3977.7778
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 11933.3334
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -16873.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1148.3334
$ws.Range("I46").Value = 747.5
$ws.Range("K46").Value = 747.5
$ws.Range("M46").Value = -559.5
$ws.Range("H122").Value = 6968
$ws.Range("I122").Value = 7484.2104
$ws.Range("J122").Value = 5333.3335
$ws.Range("K122").Value = 22452.6312
$ws.Range("L122").Value = 16000.0005
$ws.Range("M122").Value = -20002.6312
$ws.Range("N122").Value = -20900.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1988.5
$ws.Range("I126").Value = 1988.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5965.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3495.5
$ws.Range("N126").Value = ""
